# Edit the "Users" sheet:
#  - C2 changes from 8 to -1
#  - Row 7 gets new data: Hanna / hanna@gmail.com / 0
#  - Rows 8-11 (previously empty placeholder rows) are removed, shrinking
#    the used range from A1:C11 to A1:C7

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Update the Result value for Katja
$ws.Range("C2").Value = -1

# Add the new row for Hanna
$ws.Range("A7").Value = "Hanna"
$ws.Range("B7").Value = "hanna@gmail.com"
$ws.Range("C7").Value = 0
# Column C carries a default "centered" style; the new row should stay
# unstyled like rows 3-6, so reset it back to the workbook default.
$ws.Range("C7").Style = "Normal"

# Remove the leftover empty rows 8-11 so the used range shrinks to A1:C7
$ws.Range("A8:C11").Delete()
